# Apply the "cleaned up my library" edit to the Pin Allocation sheet.
#
# Before (rows 20-25 on Sheet1):
#   20: more todo:
#   21: silkscreen
#   22: rename board
#   23: test all footprints
#   24: remove makerbot from botstep23 socket silk
#   25: add net classes to high voltage lines (spindle, mosfets, etc)
#
# After:
#   20: more todo:
#   21: rename board
#   22: test all footprints
#   23: add net classes to high voltage lines (spindle, mosfets, etc)
#   24: silkscreen
#
# i.e. the "remove makerbot ..." todo item is done and removed entirely,
# and the "silkscreen" item is pushed to the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the completed "remove makerbot from botstep23 socket silk" row.
# This shifts the remaining rows (25 -> 24) up by one.
$ws.Rows("24").Delete()

# Re-order the remaining todo items so "silkscreen" moves to the end of
# the list.
$ws.Range("A21").Value = "rename board"
$ws.Range("A22").Value = "test all footprints"
$ws.Range("A23").Value = "add net classes to high voltage lines (spindle, mosfets, etc)"
$ws.Range("A24").Value = "silkscreen"

# Keep the current selection in sync with where Excel would have left the
# cursor after these edits.
$ws.Range("B24").Select()
